# Rename the "Plant_species" factor label to "Host species" in both ANOVA
# tables on the Shannon_ANOVA sheet, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shannon_ANOVA")

# Row 4 (Shannon Index table) and row 10 (Inverse Simpson table) both carry
# the "Plant_species" factor label in column A - rename to "Host species".
$ws.Range("A4").Value = "Host species"
$ws.Range("A10").Value = "Host species"

# Update the selected cell shown when the workbook is next opened.
[void]$ws.Range("K4").Select()
